# Update the kmeans process_results sheet to reflect results for a new
# test image ("test_image_1.tif") - successfully read the scale bar on
# one test image.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fname column (A2:A5) now points at the new test image ---
$ws.Range("A2").Value = "test_image_1.tif"
$ws.Range("A3").Value = "test_image_1.tif"
$ws.Range("A4").Value = "test_image_1.tif"
$ws.Range("A5").Value = "test_image_1.tif"

# --- Row 2 ---
$ws.Range("D2").Value = 0.80971659919028338
$ws.Range("E2").Value = 655.0607287449393
$ws.Range("F2").Value = 549.79757085020242
$ws.Range("G2").Value = 1.1914580265095729
$ws.Range("H2").Value = 192037
$ws.Range("I2").Value = 400.38754395112579
$ws.Range("J2").Value = 125907.32514874852
$ws.Range("K2").Value = 214.02914011958359
$ws.Range("L2").Value = 5261.538461538461
$ws.Range("M2").Value = 0.057152517171741753
$ws.Range("N2").Value = 0.078196291144331062
$ws.Range("O2").Value = 334.60599780250681
$ws.Range("P2").Value = 549.65203059827013
$ws.Range("Q2").Value = 25.50559528712196
$ws.Range("R2").Value = 25.50559528712196

# --- Row 3 ---
$ws.Range("D3").Value = 0.80971659919028338
$ws.Range("E3").Value = 295.54655870445345
$ws.Range("F3").Value = 225.91093117408906
$ws.Range("G3").Value = 1.3082437275985663
$ws.Range("H3").Value = 49618
$ws.Range("I3").Value = 203.52029763931074
$ws.Range("J3").Value = 32531.593699290264
$ws.Range("K3").Value = 82.333224198312294
$ws.Range("L3").Value = 1410.5263157894738
$ws.Range("M3").Value = 0.20547208005607173
$ws.Range("N3").Value = 0.34600305427999284
$ws.Range("O3").Value = 839.22759885525409
$ws.Range("P3").Value = 468.38020476439999
$ws.Range("Q3").Value = 30.221590338439693
$ws.Range("R3").Value = 30.221590338439693

# --- Row 4 ---
$ws.Range("D4").Value = 0.80971659919028338
$ws.Range("E4").Value = 9.7165991902834001
$ws.Range("F4").Value = 6.4777327935222671
$ws.Range("G4").Value = 1.5
$ws.Range("H4").Value = 110
$ws.Range("I4").Value = 9.5826239234059862
$ws.Range("J4").Value = 72.120506810470573
$ws.Range("K4").Value = 3.5330940923602463
$ws.Range("L4").Value = 31.247639500791408
$ws.Range("M4").Value = 0.92818426720588398
$ws.Range("N4").Value = 0.18568268190298765
$ws.Range("O4").Value = 317.92727272727274
$ws.Range("P4").Value = 616.25454545454545
$ws.Range("Q4").Value = 6.851449106449107
$ws.Range("R4").Value = 6.851449106449107

# --- Row 5 ---
$ws.Range("D5").Value = 0.80971659919028338
$ws.Range("E5").Value = 10.526315789473683
$ws.Range("F5").Value = 8.9068825910931171
$ws.Range("G5").Value = 1.1818181818181817
$ws.Range("H5").Value = 107
$ws.Range("I5").Value = 9.4510484721611245
$ws.Range("J5").Value = 70.15358389745775
$ws.Range("K5").Value = 3.7963541089374435
$ws.Range("L5").Value = 34.355278485492676
$ws.Range("M5").Value = 0.74691774936711841
$ws.Range("N5").Value = 0.20968532818614247
$ws.Range("O5").Value = 550.8878504672897
$ws.Range("P5").Value = 629.12149532710282
$ws.Range("Q5").Value = 3.3691803278688521
$ws.Range("R5").Value = 3.3691803278688521

# --- Column widths re-sized to fit the new data ---
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 2.0
$ws.Columns.Item(3).ColumnWidth = 6.333333333333333
$ws.Columns.Item(4).ColumnWidth = 11.833333333333334
$ws.Columns.Item(5).ColumnWidth = 10.833333333333334
$ws.Columns.Item(6).ColumnWidth = 10.833333333333334
$ws.Columns.Item(7).ColumnWidth = 11.166666666666666
$ws.Columns.Item(8).ColumnWidth = 10.666666666666666
$ws.Columns.Item(9).ColumnWidth = 10.833333333333334
$ws.Columns.Item(10).ColumnWidth = 10.833333333333334
$ws.Columns.Item(11).ColumnWidth = 10.833333333333334
$ws.Columns.Item(12).ColumnWidth = 10.833333333333334
$ws.Columns.Item(13).ColumnWidth = 12.833333333333334
$ws.Columns.Item(14).ColumnWidth = 12.833333333333334
$ws.Columns.Item(15).ColumnWidth = 13.5
$ws.Columns.Item(16).ColumnWidth = 13.5
$ws.Columns.Item(17).ColumnWidth = 10.833333333333334
$ws.Columns.Item(18).ColumnWidth = 10.833333333333334
